$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$High = "High"
$Medium = "Medium"
$Low = "Low"

# Row 3: Data Import
$ws.Range("B3").Value = $High
$ws.Range("C3").Value = $High
$ws.Range("D3").Value = $High

# Row 4: Data Manipulation
$ws.Range("B4").Value = $High
$ws.Range("C4").Value = $High
$ws.Range("D4").Value = $High

# Row 5: Data Visualization
$ws.Range("B5").Value = $High
$ws.Range("C5").Value = $High
$ws.Range("D5").Value = $High

# Row 6: Data Reporting
$ws.Range("B6").Value = $Medium
$ws.Range("C6").Value = $High
$ws.Range("D6").Value = $Low

# Row 7: Basic Modeling
$ws.Range("B7").Value = $High
$ws.Range("C7").Value = $Low
$ws.Range("D7").Value = $High

# Row 8: Advanced Techniques
$ws.Range("B8").Value = $High
$ws.Range("C8").Value = $Low
$ws.Range("D8").Value = $High

# Row 10: Flexibility
$ws.Range("B10").Value = $Medium
$ws.Range("C10").Value = $Low
$ws.Range("D10").Value = $High

# Row 11: Reusability
$ws.Range("B11").Value = $Medium
$ws.Range("C11").Value = $Low
$ws.Range("D11").Value = $High

# Update selection to match the target view state
$ws.Range("B3").Select()
